$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1465.0834
$ws.Range("I40").Value = 1414.7273
$ws.Range("J40").Value = 1507.6923
$ws.Range("K40").Value = 1414.7273
$ws.Range("L40").Value = 1507.6923
$ws.Range("M40").Value = -1239.7273
$ws.Range("N40").Value = -1857.6923

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 15386946
$ws.Range("I116").Value = 40001760
$ws.Range("J116").Value = 2686.875
$ws.Range("K116").Value = 40001760
$ws.Range("L116").Value = 2686.875
$ws.Range("M116").Value = -39998318
$ws.Range("N116").Value = -9570.875

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2668.2163
$ws.Range("I132").Value = 2136.3635
$ws.Range("K132").Value = 6409.0905
$ws.Range("M132").Value = -3879.0905

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1536
$ws.Range("I137").Value = 1572.2916
$ws.Range("J137").Value = 1492.45
$ws.Range("K137").Value = 4716.8748
$ws.Range("L137").Value = 4477.35
$ws.Range("M137").Value = -2166.8748
$ws.Range("N137").Value = -9577.35

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7595.396
$ws.Range("I32").Value = 7997.8623
$ws.Range("J32").Value = 5583.0625
$ws.Range("K32").Value = 7997.8623
$ws.Range("L32").Value = 5583.0625
$ws.Range("M32").Value = -7710.8623
$ws.Range("N32").Value = -6157.0625

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1144.8182
$ws.Range("I61").Value = 1167.4286
$ws.Range("J61").Value = 1105.25
$ws.Range("K61").Value = 1167.4286
$ws.Range("L61").Value = 1105.25
$ws.Range("M61").Value = -955.4286
$ws.Range("N61").Value = -1529.25

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4388.125
$ws.Range("I63").Value = 5061
$ws.Range("J63").Value = 3266.6667
$ws.Range("K63").Value = 5061
$ws.Range("L63").Value = 3266.6667
$ws.Range("M63").Value = -4375
$ws.Range("N63").Value = -4638.6667

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4388.125
$ws.Range("I66").Value = 5061
$ws.Range("J66").Value = 3266.6667
$ws.Range("K66").Value = 25305
$ws.Range("L66").Value = 16333.3335
$ws.Range("M66").Value = -21873
$ws.Range("N66").Value = -23197.3335

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 933.1818
$ws.Range("I74").Value = 896.25
$ws.Range("J74").Value = 1140
$ws.Range("K74").Value = 896.25
$ws.Range("L74").Value = 1140
$ws.Range("M74").Value = -22.25
$ws.Range("N74").Value = -2888

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 933.1818
$ws.Range("I77").Value = 896.25
$ws.Range("J77").Value = 1140
$ws.Range("K77").Value = 4481.25
$ws.Range("L77").Value = 5700
$ws.Range("M77").Value = -113.25
$ws.Range("N77").Value = -14436

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2418.6155
$ws.Range("I132").Value = 2364
$ws.Range("J132").Value = 2447.5293
$ws.Range("K132").Value = 7092
$ws.Range("L132").Value = 7342.5879
$ws.Range("M132").Value = -4562
$ws.Range("N132").Value = -12402.5879

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1144.8182
$ws.Range("I136").Value = 1167.4286
$ws.Range("J136").Value = 1105.25
$ws.Range("K136").Value = 3502.2858
$ws.Range("L136").Value = 3315.75
$ws.Range("M136").Value = -952.2857999999997
$ws.Range("N136").Value = -8415.75

# BSM row 95
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 26500
$ws.Range("J95").Value = 26500
$ws.Range("L95").Value = 26500
$ws.Range("N95").Value = -31992

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2125
$ws.Range("I134").Value = 1872.0938
$ws.Range("J134").Value = 3024.2222
$ws.Range("K134").Value = 5616.2814
$ws.Range("L134").Value = 9072.6666
$ws.Range("M134").Value = -3081.2814
$ws.Range("N134").Value = -14142.6666

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1752.75
$ws.Range("I31").Value = 1258.375
$ws.Range("J31").Value = 2741.5
$ws.Range("K31").Value = 1258.375
$ws.Range("L31").Value = 2741.5
$ws.Range("M31").Value = -963.375
$ws.Range("N31").Value = -3331.5

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1752.75
$ws.Range("I34").Value = 1258.375
$ws.Range("J34").Value = 2741.5
$ws.Range("K34").Value = 1258.375
$ws.Range("L34").Value = 2741.5
$ws.Range("M34").Value = -1056.375
$ws.Range("N34").Value = -3145.5

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1030275.25
$ws.Range("I58").Value = 1611345.6
$ws.Range("J58").Value = 2227.8462
$ws.Range("K58").Value = 1611345.6
$ws.Range("L58").Value = 2227.8462
$ws.Range("M58").Value = -1611142.6
$ws.Range("N58").Value = -2633.8462

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 411076.1
$ws.Range("I132").Value = 521024.62
$ws.Range("J132").Value = 2695.8572
$ws.Range("K132").Value = 1563073.86
$ws.Range("L132").Value = 8087.571599999999
$ws.Range("M132").Value = -1560543.86
$ws.Range("N132").Value = -13147.5716

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1480.68
$ws.Range("I134").Value = 1209.1945
$ws.Range("J134").Value = 2178.7856
$ws.Range("K134").Value = 3627.5835
$ws.Range("L134").Value = 6536.3568
$ws.Range("M134").Value = -1092.5835
$ws.Range("N134").Value = -11606.3568

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1030275.25
$ws.Range("I136").Value = 1611345.6
$ws.Range("J136").Value = 2227.8462
$ws.Range("K136").Value = 4834036.800000001
$ws.Range("L136").Value = 6683.5386
$ws.Range("M136").Value = -4831486.800000001
$ws.Range("N136").Value = -11783.5386

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 709.6316
$ws.Range("I122").Value = 500.22223
$ws.Range("J122").Value = 898.1
$ws.Range("K122").Value = 4502.00007
$ws.Range("L122").Value = 8082.900000000001
$ws.Range("M122").Value = -2052.00007
$ws.Range("N122").Value = -12982.9

# GSM row 39
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 12500
$ws.Range("J39").Value = 12500
$ws.Range("L39").Value = 12500
$ws.Range("N39").Value = -13564

# GSM row 95
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 23914.666
$ws.Range("J95").Value = 23914.666
$ws.Range("L95").Value = 23914.666
$ws.Range("N95").Value = -29406.666

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 684.3333
$ws.Range("I93").Value = 684.3333
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 684.3333
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 563.6667
$ws.Range("N93").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5695.1113
$ws.Range("I132").Value = 5700.4
$ws.Range("K132").Value = 17101.2
$ws.Range("M132").Value = -14571.2

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 23946.75
$ws.Range("J140").Value = 23946.75
$ws.Range("L140").Value = 23946.75
$ws.Range("N140").Value = -34306.75

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1145.6818
$ws.Range("I136").Value = 1045.3
$ws.Range("J136").Value = 2149.5
$ws.Range("K136").Value = 3135.9
$ws.Range("L136").Value = 6448.5
$ws.Range("M136").Value = -585.8999999999996
$ws.Range("N136").Value = -11548.5

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 69313.336
$ws.Range("J139").Value = 69313.336
$ws.Range("L139").Value = 69313.336
$ws.Range("N139").Value = -79593.336
